$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at the top of the data block (row 6) ---
# This is a new weekly report row; everything that was previously in
# rows 6-17 shifts down to rows 7-18.
$ws.Rows("6:6").Insert()

$ws.Range("A6").Value = 9
$ws.Range("B6").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C6").Value = "Metropolitana"
$ws.Range("D6").Value = 44467
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 100112010
$ws.Range("G6").Value = "Achicoria"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 52
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5500
$ws.Range("N6").Value = "$/caja 16 unidades"
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 344
$ws.Range("Q6").Value = 16
$ws.Range("R6").Value = "Hortaliza"

# --- Insert another new row before the last existing data row ---
# The original last row (previously row 17, now row 18 after the first
# insert) moves down to row 19; the new row is inserted at row 18.
$ws.Rows("18:18").Insert()

$ws.Range("A18").Value = 9
$ws.Range("B18").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C18").Value = "Metropolitana"
$ws.Range("D18").Value = 44376
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 100112010
$ws.Range("G18").Value = "Achicoria"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 43
$ws.Range("K18").Value = 4500
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 4756
$ws.Range("N18").Value = "$/caja 16 unidades"
$ws.Range("O18").Value = "Provincia de Quillota"
$ws.Range("P18").Value = 297
$ws.Range("Q18").Value = 16
$ws.Range("R18").Value = "Hortaliza"
